$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.673.26'
$ws.Range("E2").Value = '  +1.56%  '
$ws.Range("D3").Value = '1.638.02'
$ws.Range("E3").Value = '  +1.94%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("E7").Value = '  +1.19%  '
$ws.Range("E8").Value = '  +1.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0623'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.33%  '
$ws.Range("E11").Value = '  +2.62%  '
$ws.Range("D12").Value = '1.866.23'
$ws.Range("E12").Value = '  +1.85%  '
$ws.Range("D13").Value = '1.634.25'
$ws.Range("E13").Value = '  +1.67%  '
$ws.Range("E14").Value = '  +0.64%  '
$ws.Range("E15").Value = '  +1.89%  '
$ws.Range("D16").Value = '26.671.00'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.03'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.52%  '
$ws.Range("D18").Value = '0.0₃0740'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '210.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.02%  '
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.98%  '
$ws.Range("E22").Value = '  +0.99%  '
$ws.Range("E23").Value = '  +2.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.02'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("E27").Value = '  -0.73%  '
$ws.Range("E28").Value = '  +2.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0517'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.83%  '
$ws.Range("E31").Value = '  -0.54%  '
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("E34").Value = '  +0.86%  '
$ws.Range("E35").Value = '  -0.95%  '
$ws.Range("D36").Value = '1.169.93'
$ws.Range("E36").Value = '  +0.56%  '
$ws.Range("E37").Value = '  +0.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.807'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("E41").Value = '  +0.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.795'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("E43").Value = '  +1.70%  '
$ws.Range("D44").Value = '1.774.38'
$ws.Range("E44").Value = '  +1.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("E46").Value = '  +2.13%  '
$ws.Range("E47").Value = '  +5.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.410'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.99%  '
